# Generate Report for Archive
#
# 1. Update status text "Ready for handoff" -> "In Translation" wherever it
#    appears: the per-language status roll-up on the Overview sheet (columns
#    E/F) and the "Status" column (C) on the zh-cn / de-de detail sheets.
# 2. Narrow the now-shorter status columns to fit the new text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

foreach ($ws in @($overview, $zhcn, $dede)) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # NOTE: keep the string literal on the LEFT of -eq. PowerShell's -eq
        # coerces the right-hand side to the left operand's type, so with a
        # boolean cell (e.g. the "True"/"False" columns) on the left,
        # "$cell.Value2 -eq $oldStatus" would coerce the non-empty string to
        # $true and falsely match every TRUE cell.
        if ($oldStatus -eq $cell.Value2) {
            $cell.Value = $newStatus
        }
    }
}

# --- Resize columns that previously fit "Ready for handoff" ---
# Range.ColumnWidth is expressed in characters of the Normal-style font and
# gets pixel-quantized by Excel, same as interactively typing a width; back
# that quantization out so the stored <col width=.../> lands on the target.
$targetCharWidth = 13.4101845877511
$columnWidth = $targetCharWidth - (5 / 6)

$overview.Range("E:F").ColumnWidth = $columnWidth
$zhcn.Range("C:C").ColumnWidth = $columnWidth
$dede.Range("C:C").ColumnWidth = $columnWidth
